# Update crypto price/volume data (and swap the BabyDogeCoin / dogwifhat rows) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '67.921.76'
    'E2' = '  +0.92%  '
    'D3' = '2.641.89'
    'E3' = '  +0.80%  '
    'D5' = '598.20'
    'E5' = '  +0.27%  '
    'D6' = '153.82'
    'E6' = '  +0.84%  '
    'E7' = '  -0.02%  '
    'D8' = '0.551'
    'E8' = '  -0.31%  '
    'D9' = '2.642.75'
    'E9' = '  +0.84%  '
    'E10' = '  +10.57%  '
    'D12' = '5.22'
    'E12' = '  +0.75%  '
    'E13' = '  -0.01%  '
    'D14' = '27.66'
    'E14' = '  +0.34%  '
    'E15' = '  +4.00%  '
    'D16' = '3.120.40'
    'E16' = '  +0.86%  '
    'D17' = '67.892.51'
    'E17' = '  +0.94%  '
    'D18' = '2.640.30'
    'E18' = '  +0.60%  '
    'D19' = '11.47'
    'E19' = '  +3.31%  '
    'D20' = '373.24'
    'E20' = '  +2.83%  '
    'D21' = '7.50'
    'E21' = '  +0.20%  '
    'D22' = '4.25'
    'E22' = '  -0.86%  '
    'E23' = '  -1.31%  '
    'D24' = '2.06'
    'E24' = '  -1.63%  '
    'D25' = '72.17'
    'E25' = '  +1.64%  '
    'E26' = '  +0.26%  '
    'D27' = '9.99'
    'E27' = '  -0.80%  '
    'E28' = '  -0.08%  '
    'E29' = '  +2.21%  '
    'E30' = '  +0.28%  '
    'D31' = '577.48'
    'E31' = '  -1.14%  '
    'E32' = '  +0.76%  '
    'D33' = '7.89'
    'E33' = '  +1.00%  '
    'E34' = '  +0.35%  '
    'E35' = '  -0.02%  '
    'D36' = '0.126'
    'E36' = '  -0.25%  '
    'E37' = '  -0.01%  '
    'D38' = '157.94'
    'E38' = '  +0.33%  '
    'D39' = '19.23'
    'E39' = '  +0.46%  '
    'E40' = '  +5.40%  '
    'D41' = '0.370'
    'E41' = '  +0.47%  '
    'D42' = '5.37'
    'E42' = '  +2.11%  '
    'B43' = 'BabyDogeCoin'
    'C43' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D43' = '0.0₆0336'
    'E43' = '  +17.24%  '
    'B44' = 'dogwifhat'
    'C44' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D44' = '2.63'
    'E44' = '  +2.36%  '
    'D45' = '17.12'
    'E45' = '  +4.73%  '
    'D46' = '1.00'
    'E46' = '  +0.07%  '
    'D47' = '40.25'
    'E47' = '  -2.29%  '
    'D48' = '156.44'
    'E48' = '  +0.10%  '
    'E49' = '  -0.69%  '
    'D50' = '22.01'
    'E50' = '  +7.08%  '
    'D51' = '1.71'
    'E51' = '  -0.99%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text so numeric-looking strings (e.g. '598.20', '1.00') are not
    # coerced to floating point values; ClearFormats() afterwards restores the
    # cell's original (default) style so only the content changes.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
